# Updated cryptos list on Sat May  4 22:41:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text storage to preserve formatting like
# thousand separators written with dots (e.g. 63.831.89) instead of Excel
# reinterpreting them as numbers.
$priceUpdates = @{
    "D2" = '63.831.89'
    "D3" = '3.115.79'
    "D5" = '585.57'
    "D6" = '146.24'
    "D8" = '3.110.57'
    "D10" = '0.161'
    "D11" = '5.75'
    "D13" = '0.0000251'
    "D14" = '37.02'
    "D16" = '3.634.74'
    "D17" = '63.742.09'
    "D18" = '7.14'
    "D19" = '3.112.38'
    "D20" = '464.49'
    "D21" = '14.31'
    "D22" = '0.730'
    "D24" = '13.11'
    "D25" = '82.00'
    "D28" = '2.69'
    "D31" = '6.86'
    "D32" = '26.96'
    "D34" = '0.0₃0871'
    "D35" = '2.36'
    "D37" = '3.42'
    "D39" = '50.93'
    "D40" = '448.77'
    "D41" = '8.68'
    "D43" = '2.881.94'
    "D46" = '2.16'
    "D47" = '35.81'
    "D49" = '123.90'
    "D51" = '24.65'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = $origStyle
}

# Volume(1h) (column E) updates - plain text percentages, already safe from
# Excel's numeric auto-detection because of the surrounding whitespace.
$volumeUpdates = @{
    "E3" = '  -0.08%  '
    "E4" = '  -0.04%  '
    "E5" = '  -0.26%  '
    "E6" = '  +1.91%  '
    "E8" = '  +0.15%  '
    "E9" = '  -0.46%  '
    "E10" = '  +10.86%  '
    "E11" = '  -0.14%  '
    "E12" = '  -0.67%  '
    "E13" = '  +2.96%  '
    "E14" = '  +4.07%  '
    "E15" = '  -0.89%  '
    "E16" = '  +0.08%  '
    "E17" = '  +1.15%  '
    "E18" = '  -1.97%  '
    "E19" = '  -0.13%  '
    "E20" = '  +2.41%  '
    "E21" = '  +1.57%  '
    "E22" = '  -0.64%  '
    "E23" = '  -0.02%  '
    "E24" = '  -3.85%  '
    "E25" = '  -0.03%  '
    "E27" = '  +8.19%  '
    "E28" = '  -0.42%  '
    "E29" = '  -1.53%  '
    "E30" = '  -0.05%  '
    "E31" = '  +0.18%  '
    "E33" = '  -3.05%  '
    "E34" = '  +8.06%  '
    "E35" = '  +2.85%  '
    "E36" = '  +0.70%  '
    "E37" = '  +12.78%  '
    "E38" = '  -0.01%  '
    "E39" = '  +0.17%  '
    "E40" = '  +4.50%  '
    "E41" = '  -1.25%  '
    "E42" = '  -0.86%  '
    "E43" = '  -2.30%  '
    "E44" = '  -0.30%  '
    "E45" = '  -0.71%  '
    "E46" = '  -0.33%  '
    "E47" = '  +3.03%  '
    "E48" = '  +0.03%  '
    "E49" = '  -1.63%  '
    "E50" = '  -0.64%  '
    "E51" = '  -0.62%  '
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

